$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 updates
$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 2.63
$ws.Range("Q5").Value = 1.8
$ws.Range("R5").Value = 2.05
$ws.Range("S5").Value = 2.35
$ws.Range("T5").Value = 1.57
$ws.Range("W5").Value = 4.5
$ws.Range("X5").Value = 1.18
$ws.Range("AF5").Value = 11

# Row 6 updates
$ws.Range("G6").Value = 2.6
$ws.Range("Q6").Value = 1.88
$ws.Range("R6").Value = 1.98
$ws.Range("U6").Value = 3.95
$ws.Range("V6").Value = 1.26
$ws.Range("AG6").Value = 23
$ws.Range("AM6").Value = 451
$ws.Range("AQ6").Value = 34

# Row 8 updates
$ws.Range("G8").Value = 3.5
$ws.Range("I8").Value = 2.3
$ws.Range("L8").Value = 3.2
$ws.Range("M8").Value = 1.13
$ws.Range("N8").Value = 6
$ws.Range("Q8").Value = 2.1
$ws.Range("R8").Value = 1.78
$ws.Range("U8").Value = 4.3
$ws.Range("V8").Value = 1.21
$ws.Range("AC8").Value = 8
$ws.Range("AL8").Value = 67
